$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "1.00", "5.92")
# must be forced to Text format first so Excel keeps them as text instead of
# auto-converting to a number (which would drop formatting like trailing zeros).
$textForceValues = @{
    'D4' = '1.00'
    'D5' = '304.22'
    'D6' = '95.11'
    'D10' = '34.70'
    'D13' = '6.75'
    'D17' = '0.786'
    'D19' = '12.22'
    'D21' = '5.92'
    'D22' = '68.31'
    'D23' = '236.23'
    'D27' = '23.44'
    'D28' = '36.22'
    'D30' = '9.37'
    'D31' = '159.81'
    'D33' = '5.17'
    'D37' = '16.81'
    'D44' = '0.0280'
    'D45' = '18.56'
    'D46' = '9.86'
    'D48' = '52.69'
    'D49' = '72.13'
    'D50' = '90.83'
    'D51' = '1.12'
}
foreach ($cellRef in $textForceValues.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForceValues[$cellRef]
}

# Remaining cells: plain text values that are safe to assign directly.
$plainValues = @{
    'D2' = '41.862.16'
    'E2' = '  -0.61%  '
    'D3' = '2.251.11'
    'E3' = '  -0.62%  '
    'E4' = '  -0.07%  '
    'E5' = '  -0.55%  '
    'E6' = '  +1.13%  '
    'E7' = '  -1.54%  '
    'E8' = '  -0.03%  '
    'E9' = '  +0.16%  '
    'E10' = '  +4.59%  '
    'E11' = '  -2.01%  '
    'E12' = '  -0.05%  '
    'E13' = '  +1.17%  '
    'D14' = '2.599.10'
    'E14' = '  -0.59%  '
    'E15' = '  +0.05%  '
    'D16' = '2.242.22'
    'E16' = '  -1.33%  '
    'E17' = '  -0.05%  '
    'D18' = '41.743.12'
    'E18' = '  -0.60%  '
    'E19' = '  -3.68%  '
    'D20' = '0.0₃0898'
    'E22' = '  +0.29%  '
    'E23' = '  -3.26%  '
    'E24' = '  -2.51%  '
    'E25' = '  +0.07%  '
    'E26' = '  -1.93%  '
    'E27' = '  -2.25%  '
    'E28' = '  +3.52%  '
    'E29' = '  +1.36%  '
    'E30' = '  -3.34%  '
    'E32' = '  -0.05%  '
    'E33' = '  -3.31%  '
    'E34' = '  +2.59%  '
    'E35' = '  -1.91%  '
    'E36' = '  +0.05%  '
    'E37' = '  -1.32%  '
    'E38' = '  -1.62%  '
    'E39' = '  -2.27%  '
    'E40' = '  -0.30%  '
    'E41' = '  -3.16%  '
    'E42' = '  +0.52%  '
    'D43' = '1.956.86'
    'E43' = '  -2.46%  '
    'E44' = '  -1.37%  '
    'E45' = '  -5.51%  '
    'E46' = '  -3.29%  '
    'E47' = '  -1.63%  '
    'E48' = '  -1.80%  '
    'E49' = '  -1.30%  '
    'E50' = '  -0.83%  '
    'B51' = 'TrustWalletToken'
    'C51' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E51' = '  -2.23%  '
}
foreach ($cellRef in $plainValues.Keys) {
    $ws.Range($cellRef).Value = $plainValues[$cellRef]
}
